$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds a new weekly price observation (row 152) for "Acelga" at
# "Feria Lagunitas de Puerto Montt" and shifts every existing observation in
# rows 39-151 down by one row (newest entry goes on top, row 39).

# --- Row 152: constant/descriptive columns, copied from the surrounding rows ---
$ws.Range("A152").Value = 4
$ws.Range("B152").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C152").Value = "Los Lagos"
$ws.Range("E152").Value = 10
$ws.Range("F152").Value = 100112009
$ws.Range("G152").Value = "Acelga"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("R152").Value = "Hortaliza"
$ws.Range("D152").NumberFormat = $ws.Range("D151").NumberFormat

# --- Data columns (D,J,K,L,M,N,O,P,Q): row 39 gets the new observation,
# rows 40-152 get the value that used to sit one row above them ---
$ws.Range("D39").Value = 44622
$ws.Range("J39").Value = 20
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = 10000
$ws.Range("N39").Value = "$/docena de atados (12 kilos)"
$ws.Range("O39").Value = "Región de La Araucanía"
$ws.Range("P39").Value = 833
$ws.Range("Q39").Value = 12

$ws.Range("D40").Value = 44235
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 800
$ws.Range("M40").Value = 800
$ws.Range("N40").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O40").Value = "Región de La Araucanía"
$ws.Range("P40").Value = 533
$ws.Range("Q40").Value = 1.5

$ws.Range("D41").Value = 44209
$ws.Range("J41").Value = 100
$ws.Range("K41").Value = 800
$ws.Range("L41").Value = 800
$ws.Range("M41").Value = 800
$ws.Range("N41").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O41").Value = "Región de La Araucanía"
$ws.Range("P41").Value = 533
$ws.Range("Q41").Value = 1.5

$ws.Range("D42").Value = 44250
$ws.Range("J42").Value = 250
$ws.Range("K42").Value = 800
$ws.Range("L42").Value = 1000
$ws.Range("M42").Value = 920
$ws.Range("N42").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O42").Value = "Región de La Araucanía"
$ws.Range("P42").Value = 613
$ws.Range("Q42").Value = 1.5

$ws.Range("D43").Value = 44253
$ws.Range("J43").Value = 300
$ws.Range("K43").Value = 800
$ws.Range("L43").Value = 800
$ws.Range("M43").Value = 800
$ws.Range("N43").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O43").Value = "Región de La Araucanía"
$ws.Range("P43").Value = 533
$ws.Range("Q43").Value = 1.5

$ws.Range("D44").Value = 44298
$ws.Range("J44").Value = 60
$ws.Range("K44").Value = 800
$ws.Range("L44").Value = 800
$ws.Range("M44").Value = 800
$ws.Range("N44").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O44").Value = "Región de La Araucanía"
$ws.Range("P44").Value = 533
$ws.Range("Q44").Value = 1.5

$ws.Range("D45").Value = 44211
$ws.Range("J45").Value = 150
$ws.Range("K45").Value = 800
$ws.Range("L45").Value = 800
$ws.Range("M45").Value = 800
$ws.Range("N45").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O45").Value = "Región de La Araucanía"
$ws.Range("P45").Value = 533
$ws.Range("Q45").Value = 1.5

$ws.Range("D46").Value = 44614
$ws.Range("J46").Value = 80
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = 10000
$ws.Range("N46").Value = "$/docena de atados (12 kilos)"
$ws.Range("O46").Value = "Región de La Araucanía"
$ws.Range("P46").Value = 833
$ws.Range("Q46").Value = 12

$ws.Range("D47").Value = 44424
$ws.Range("J47").Value = 50
$ws.Range("K47").Value = 4000
$ws.Range("L47").Value = 4000
$ws.Range("M47").Value = 4000
$ws.Range("N47").Value = "$/docena de atados (4 kilos)"
$ws.Range("O47").Value = "Región del Maule"
$ws.Range("P47").Value = 1000
$ws.Range("Q47").Value = 4

$ws.Range("D48").Value = 44495
$ws.Range("J48").Value = 200
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 3000
$ws.Range("M48").Value = 3000
$ws.Range("N48").Value = "$/docena de atados (4 kilos)"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 750
$ws.Range("Q48").Value = 4

$ws.Range("D49").Value = 44306
$ws.Range("J49").Value = 300
$ws.Range("K49").Value = 800
$ws.Range("L49").Value = 800
$ws.Range("M49").Value = 800
$ws.Range("N49").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O49").Value = "Región de La Araucanía"
$ws.Range("P49").Value = 533
$ws.Range("Q49").Value = 1.5

$ws.Range("D50").Value = 44222
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 700
$ws.Range("L50").Value = 800
$ws.Range("M50").Value = 750
$ws.Range("N50").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O50").Value = "Región de La Araucanía"
$ws.Range("P50").Value = 500
$ws.Range("Q50").Value = 1.5

$ws.Range("D51").Value = 44383
$ws.Range("J51").Value = 200
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 3500
$ws.Range("M51").Value = 3500
$ws.Range("N51").Value = "$/docena de atados (4 kilos)"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 875
$ws.Range("Q51").Value = 4

$ws.Range("D52").Value = 44469
$ws.Range("J52").Value = 100
$ws.Range("K52").Value = 4000
$ws.Range("L52").Value = 4000
$ws.Range("M52").Value = 4000
$ws.Range("N52").Value = "$/docena de atados (4 kilos)"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 1000
$ws.Range("Q52").Value = 4

$ws.Range("D53").Value = 44295
$ws.Range("J53").Value = 300
$ws.Range("K53").Value = 800
$ws.Range("L53").Value = 800
$ws.Range("M53").Value = 800
$ws.Range("N53").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O53").Value = "Región de La Araucanía"
$ws.Range("P53").Value = 533
$ws.Range("Q53").Value = 1.5

$ws.Range("D54").Value = 44369
$ws.Range("J54").Value = 200
$ws.Range("K54").Value = 3500
$ws.Range("L54").Value = 3500
$ws.Range("M54").Value = 3500
$ws.Range("N54").Value = "$/docena de atados (4 kilos)"
$ws.Range("O54").Value = "Región del Maule"
$ws.Range("P54").Value = 875
$ws.Range("Q54").Value = 4

$ws.Range("D55").Value = 44278
$ws.Range("J55").Value = 250
$ws.Range("K55").Value = 800
$ws.Range("L55").Value = 800
$ws.Range("M55").Value = 800
$ws.Range("N55").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O55").Value = "Región de La Araucanía"
$ws.Range("P55").Value = 533
$ws.Range("Q55").Value = 1.5

$ws.Range("D56").Value = 44435
$ws.Range("J56").Value = 400
$ws.Range("K56").Value = 4000
$ws.Range("L56").Value = 4000
$ws.Range("M56").Value = 4000
$ws.Range("N56").Value = "$/docena de atados (4 kilos)"
$ws.Range("O56").Value = "Región del Maule"
$ws.Range("P56").Value = 1000
$ws.Range("Q56").Value = 4

$ws.Range("D57").Value = 44292
$ws.Range("J57").Value = 300
$ws.Range("K57").Value = 800
$ws.Range("L57").Value = 800
$ws.Range("M57").Value = 800
$ws.Range("N57").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O57").Value = "Región de La Araucanía"
$ws.Range("P57").Value = 533
$ws.Range("Q57").Value = 1.5

$ws.Range("D58").Value = 44314
$ws.Range("J58").Value = 60
$ws.Range("K58").Value = 800
$ws.Range("L58").Value = 800
$ws.Range("M58").Value = 800
$ws.Range("N58").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O58").Value = "Región de La Araucanía"
$ws.Range("P58").Value = 533
$ws.Range("Q58").Value = 1.5

$ws.Range("D59").Value = 44176
$ws.Range("J59").Value = 200
$ws.Range("K59").Value = 3500
$ws.Range("L59").Value = 3500
$ws.Range("M59").Value = 3500
$ws.Range("N59").Value = "$/docena de atados (4 kilos)"
$ws.Range("O59").Value = "Región del Maule"
$ws.Range("P59").Value = 875
$ws.Range("Q59").Value = 4

$ws.Range("D60").Value = 44260
$ws.Range("J60").Value = 200
$ws.Range("K60").Value = 800
$ws.Range("L60").Value = 800
$ws.Range("M60").Value = 800
$ws.Range("N60").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O60").Value = "Región de La Araucanía"
$ws.Range("P60").Value = 533
$ws.Range("Q60").Value = 1.5

$ws.Range("D61").Value = 44467
$ws.Range("J61").Value = 200
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = 4000
$ws.Range("N61").Value = "$/docena de atados (4 kilos)"
$ws.Range("O61").Value = "Región del Maule"
$ws.Range("P61").Value = 1000
$ws.Range("Q61").Value = 4

$ws.Range("D62").Value = 44203
$ws.Range("J62").Value = 150
$ws.Range("K62").Value = 800
$ws.Range("L62").Value = 800
$ws.Range("M62").Value = 800
$ws.Range("N62").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O62").Value = "Región de La Araucanía"
$ws.Range("P62").Value = 533
$ws.Range("Q62").Value = 1.5

$ws.Range("D63").Value = 44341
$ws.Range("J63").Value = 250
$ws.Range("K63").Value = 800
$ws.Range("L63").Value = 800
$ws.Range("M63").Value = 800
$ws.Range("N63").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O63").Value = "Región de La Araucanía"
$ws.Range("P63").Value = 533
$ws.Range("Q63").Value = 1.5

$ws.Range("D64").Value = 44245
$ws.Range("J64").Value = 150
$ws.Range("K64").Value = 800
$ws.Range("L64").Value = 800
$ws.Range("M64").Value = 800
$ws.Range("N64").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O64").Value = "Región de La Araucanía"
$ws.Range("P64").Value = 533
$ws.Range("Q64").Value = 1.5

$ws.Range("D65").Value = 44305
$ws.Range("J65").Value = 72
$ws.Range("K65").Value = 800
$ws.Range("L65").Value = 800
$ws.Range("M65").Value = 800
$ws.Range("N65").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O65").Value = "Región de La Araucanía"
$ws.Range("P65").Value = 533
$ws.Range("Q65").Value = 1.5

$ws.Range("D66").Value = 44532
$ws.Range("J66").Value = 50
$ws.Range("K66").Value = 3500
$ws.Range("L66").Value = 3500
$ws.Range("M66").Value = 3500
$ws.Range("N66").Value = "$/docena de atados (4 kilos)"
$ws.Range("O66").Value = "Región del Maule"
$ws.Range("P66").Value = 875
$ws.Range("Q66").Value = 4

$ws.Range("D67").Value = 44540
$ws.Range("J67").Value = 200
$ws.Range("K67").Value = 3500
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = 3500
$ws.Range("N67").Value = "$/docena de atados (4 kilos)"
$ws.Range("O67").Value = "Región del Maule"
$ws.Range("P67").Value = 875
$ws.Range("Q67").Value = 4

$ws.Range("D68").Value = 44474
$ws.Range("J68").Value = 200
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = 4000
$ws.Range("N68").Value = "$/docena de atados (4 kilos)"
$ws.Range("O68").Value = "Región del Maule"
$ws.Range("P68").Value = 1000
$ws.Range("Q68").Value = 4

$ws.Range("D69").Value = 44574
$ws.Range("J69").Value = 40
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("M69").Value = 10000
$ws.Range("N69").Value = "$/docena de atados (12 kilos)"
$ws.Range("O69").Value = "Región de La Araucanía"
$ws.Range("P69").Value = 833
$ws.Range("Q69").Value = 12

$ws.Range("D70").Value = 44246
$ws.Range("J70").Value = 300
$ws.Range("K70").Value = 800
$ws.Range("L70").Value = 1000
$ws.Range("M70").Value = 900
$ws.Range("N70").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O70").Value = "Región de La Araucanía"
$ws.Range("P70").Value = 600
$ws.Range("Q70").Value = 1.5

$ws.Range("D71").Value = 44323
$ws.Range("J71").Value = 300
$ws.Range("K71").Value = 800
$ws.Range("L71").Value = 800
$ws.Range("M71").Value = 800
$ws.Range("N71").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O71").Value = "Región de La Araucanía"
$ws.Range("P71").Value = 533
$ws.Range("Q71").Value = 1.5

$ws.Range("D72").Value = 44392
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 1200
$ws.Range("L72").Value = 1200
$ws.Range("M72").Value = 1200
$ws.Range("N72").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O72").Value = "Región Metropolitana"
$ws.Range("P72").Value = 800
$ws.Range("Q72").Value = 1.5

$ws.Range("D73").Value = 44392
$ws.Range("J73").Value = 200
$ws.Range("K73").Value = 3500
$ws.Range("L73").Value = 3500
$ws.Range("M73").Value = 3500
$ws.Range("N73").Value = "$/docena de atados (4 kilos)"
$ws.Range("O73").Value = "Región del Maule"
$ws.Range("P73").Value = 875
$ws.Range("Q73").Value = 4

$ws.Range("D74").Value = 44601
$ws.Range("J74").Value = 20
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = 10000
$ws.Range("N74").Value = "$/docena de atados (12 kilos)"
$ws.Range("O74").Value = "Región de La Araucanía"
$ws.Range("P74").Value = 833
$ws.Range("Q74").Value = 12

$ws.Range("D75").Value = 44509
$ws.Range("J75").Value = 200
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 3000
$ws.Range("M75").Value = 3000
$ws.Range("N75").Value = "$/docena de atados (4 kilos)"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 750
$ws.Range("Q75").Value = 4

$ws.Range("D76").Value = 44230
$ws.Range("J76").Value = 60
$ws.Range("K76").Value = 800
$ws.Range("L76").Value = 800
$ws.Range("M76").Value = 800
$ws.Range("N76").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O76").Value = "Provincia de Cautín"
$ws.Range("P76").Value = 533
$ws.Range("Q76").Value = 1.5

$ws.Range("D77").Value = 44316
$ws.Range("J77").Value = 240
$ws.Range("K77").Value = 800
$ws.Range("L77").Value = 800
$ws.Range("M77").Value = 800
$ws.Range("N77").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O77").Value = "Región de La Araucanía"
$ws.Range("P77").Value = 533
$ws.Range("Q77").Value = 1.5

$ws.Range("D78").Value = 44159
$ws.Range("J78").Value = 250
$ws.Range("K78").Value = 3500
$ws.Range("L78").Value = 3500
$ws.Range("M78").Value = 3500
$ws.Range("N78").Value = "$/docena de atados (4 kilos)"
$ws.Range("O78").Value = "Región del Maule"
$ws.Range("P78").Value = 875
$ws.Range("Q78").Value = 4

$ws.Range("D79").Value = 44386
$ws.Range("J79").Value = 200
$ws.Range("K79").Value = 3500
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = 3500
$ws.Range("N79").Value = "$/docena de atados (4 kilos)"
$ws.Range("O79").Value = "Región del Maule"
$ws.Range("P79").Value = 875
$ws.Range("Q79").Value = 4

$ws.Range("D80").Value = 44466
$ws.Range("J80").Value = 100
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = 4000
$ws.Range("N80").Value = "$/docena de atados (4 kilos)"
$ws.Range("O80").Value = "Región del Maule"
$ws.Range("P80").Value = 1000
$ws.Range("Q80").Value = 4

$ws.Range("D81").Value = 44322
$ws.Range("J81").Value = 120
$ws.Range("K81").Value = 800
$ws.Range("L81").Value = 800
$ws.Range("M81").Value = 800
$ws.Range("N81").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O81").Value = "Región de La Araucanía"
$ws.Range("P81").Value = 533
$ws.Range("Q81").Value = 1.5

$ws.Range("D82").Value = 44320
$ws.Range("J82").Value = 300
$ws.Range("K82").Value = 800
$ws.Range("L82").Value = 800
$ws.Range("M82").Value = 800
$ws.Range("N82").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O82").Value = "Región de La Araucanía"
$ws.Range("P82").Value = 533
$ws.Range("Q82").Value = 1.5

$ws.Range("D83").Value = 44204
$ws.Range("J83").Value = 120
$ws.Range("K83").Value = 800
$ws.Range("L83").Value = 800
$ws.Range("M83").Value = 800
$ws.Range("N83").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O83").Value = "Región de La Araucanía"
$ws.Range("P83").Value = 533
$ws.Range("Q83").Value = 1.5

$ws.Range("D84").Value = 44358
$ws.Range("J84").Value = 200
$ws.Range("K84").Value = 3500
$ws.Range("L84").Value = 3500
$ws.Range("M84").Value = 3500
$ws.Range("N84").Value = "$/docena de atados (4 kilos)"
$ws.Range("O84").Value = "Región del Maule"
$ws.Range("P84").Value = 875
$ws.Range("Q84").Value = 4

$ws.Range("D85").Value = 44313
$ws.Range("J85").Value = 300
$ws.Range("K85").Value = 800
$ws.Range("L85").Value = 800
$ws.Range("M85").Value = 800
$ws.Range("N85").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O85").Value = "Región de La Araucanía"
$ws.Range("P85").Value = 533
$ws.Range("Q85").Value = 1.5

$ws.Range("D86").Value = 44333
$ws.Range("J86").Value = 120
$ws.Range("K86").Value = 800
$ws.Range("L86").Value = 800
$ws.Range("M86").Value = 800
$ws.Range("N86").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O86").Value = "Región de La Araucanía"
$ws.Range("P86").Value = 533
$ws.Range("Q86").Value = 1.5

$ws.Range("D87").Value = 44302
$ws.Range("J87").Value = 200
$ws.Range("K87").Value = 750
$ws.Range("L87").Value = 800
$ws.Range("M87").Value = 775
$ws.Range("N87").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O87").Value = "Región de La Araucanía"
$ws.Range("P87").Value = 517
$ws.Range("Q87").Value = 1.5

$ws.Range("D88").Value = 44210
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 800
$ws.Range("M88").Value = 800
$ws.Range("N88").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O88").Value = "Región de La Araucanía"
$ws.Range("P88").Value = 533
$ws.Range("Q88").Value = 1.5

$ws.Range("D89").Value = 44291
$ws.Range("J89").Value = 150
$ws.Range("K89").Value = 800
$ws.Range("L89").Value = 800
$ws.Range("M89").Value = 800
$ws.Range("N89").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O89").Value = "Región de La Araucanía"
$ws.Range("P89").Value = 533
$ws.Range("Q89").Value = 1.5

$ws.Range("D90").Value = 44217
$ws.Range("J90").Value = 100
$ws.Range("K90").Value = 800
$ws.Range("L90").Value = 800
$ws.Range("M90").Value = 800
$ws.Range("N90").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O90").Value = "Región de La Araucanía"
$ws.Range("P90").Value = 533
$ws.Range("Q90").Value = 1.5

$ws.Range("D91").Value = 44242
$ws.Range("J91").Value = 70
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 800
$ws.Range("M91").Value = 800
$ws.Range("N91").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O91").Value = "Región de La Araucanía"
$ws.Range("P91").Value = 533
$ws.Range("Q91").Value = 1.5

$ws.Range("D92").Value = 44166
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 3500
$ws.Range("L92").Value = 3500
$ws.Range("M92").Value = 3500
$ws.Range("N92").Value = "$/docena de atados (4 kilos)"
$ws.Range("O92").Value = "Región del Maule"
$ws.Range("P92").Value = 875
$ws.Range("Q92").Value = 4

$ws.Range("D93").Value = 44348
$ws.Range("J93").Value = 300
$ws.Range("K93").Value = 800
$ws.Range("L93").Value = 800
$ws.Range("M93").Value = 800
$ws.Range("N93").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O93").Value = "Región de La Araucanía"
$ws.Range("P93").Value = 533
$ws.Range("Q93").Value = 1.5

$ws.Range("D94").Value = 44175
$ws.Range("J94").Value = 100
$ws.Range("K94").Value = 3500
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = 3500
$ws.Range("N94").Value = "$/docena de atados (4 kilos)"
$ws.Range("O94").Value = "Región del Maule"
$ws.Range("P94").Value = 875
$ws.Range("Q94").Value = 4

$ws.Range("D95").Value = 44579
$ws.Range("J95").Value = 80
$ws.Range("K95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("M95").Value = 10000
$ws.Range("N95").Value = "$/docena de atados (12 kilos)"
$ws.Range("O95").Value = "Región de La Araucanía"
$ws.Range("P95").Value = 833
$ws.Range("Q95").Value = 12

$ws.Range("D96").Value = 44515
$ws.Range("J96").Value = 100
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = 3000
$ws.Range("N96").Value = "$/docena de atados (4 kilos)"
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 750
$ws.Range("Q96").Value = 4

$ws.Range("D97").Value = 44249
$ws.Range("J97").Value = 80
$ws.Range("K97").Value = 800
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 900
$ws.Range("N97").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O97").Value = "Región de La Araucanía"
$ws.Range("P97").Value = 600
$ws.Range("Q97").Value = 1.5

$ws.Range("D98").Value = 44566
$ws.Range("J98").Value = 20
$ws.Range("K98").Value = 10000
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = 10000
$ws.Range("N98").Value = "$/docena de atados (12 kilos)"
$ws.Range("O98").Value = "Región de La Araucanía"
$ws.Range("P98").Value = 833
$ws.Range("Q98").Value = 12

$ws.Range("D99").Value = 44300
$ws.Range("J99").Value = 48
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 800
$ws.Range("M99").Value = 800
$ws.Range("N99").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O99").Value = "Región de La Araucanía"
$ws.Range("P99").Value = 533
$ws.Range("Q99").Value = 1.5

$ws.Range("D100").Value = 44578
$ws.Range("J100").Value = 40
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = 10000
$ws.Range("N100").Value = "$/docena de atados (12 kilos)"
$ws.Range("O100").Value = "Región de La Araucanía"
$ws.Range("P100").Value = 833
$ws.Range("Q100").Value = 12

$ws.Range("D101").Value = 44225
$ws.Range("J101").Value = 250
$ws.Range("K101").Value = 800
$ws.Range("L101").Value = 800
$ws.Range("M101").Value = 800
$ws.Range("N101").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O101").Value = "Región de La Araucanía"
$ws.Range("P101").Value = 533
$ws.Range("Q101").Value = 1.5

$ws.Range("D102").Value = 44411
$ws.Range("J102").Value = 200
$ws.Range("K102").Value = 4000
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = 4000
$ws.Range("N102").Value = "$/docena de atados (4 kilos)"
$ws.Range("O102").Value = "Región del Maule"
$ws.Range("P102").Value = 1000
$ws.Range("Q102").Value = 4

$ws.Range("D103").Value = 44281
$ws.Range("J103").Value = 60
$ws.Range("K103").Value = 800
$ws.Range("L103").Value = 800
$ws.Range("M103").Value = 800
$ws.Range("N103").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O103").Value = "Región de La Araucanía"
$ws.Range("P103").Value = 533
$ws.Range("Q103").Value = 1.5

$ws.Range("D104").Value = 44252
$ws.Range("J104").Value = 150
$ws.Range("K104").Value = 800
$ws.Range("L104").Value = 800
$ws.Range("M104").Value = 800
$ws.Range("N104").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O104").Value = "Región de La Araucanía"
$ws.Range("P104").Value = 533
$ws.Range("Q104").Value = 1.5

$ws.Range("D105").Value = 44271
$ws.Range("J105").Value = 200
$ws.Range("K105").Value = 800
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 900
$ws.Range("N105").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O105").Value = "Región de La Araucanía"
$ws.Range("P105").Value = 600
$ws.Range("Q105").Value = 1.5

$ws.Range("D106").Value = 44162
$ws.Range("J106").Value = 200
$ws.Range("K106").Value = 3500
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = 3500
$ws.Range("N106").Value = "$/docena de atados (4 kilos)"
$ws.Range("O106").Value = "Región del Maule"
$ws.Range("P106").Value = 875
$ws.Range("Q106").Value = 4

$ws.Range("D107").Value = 44516
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 3500
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = 3500
$ws.Range("N107").Value = "$/docena de atados (4 kilos)"
$ws.Range("O107").Value = "Región del Maule"
$ws.Range("P107").Value = 875
$ws.Range("Q107").Value = 4

$ws.Range("D108").Value = 44568
$ws.Range("J108").Value = 40
$ws.Range("K108").Value = 10000
$ws.Range("L108").Value = 10000
$ws.Range("M108").Value = 10000
$ws.Range("N108").Value = "$/docena de atados (12 kilos)"
$ws.Range("O108").Value = "Región de La Araucanía"
$ws.Range("P108").Value = 833
$ws.Range("Q108").Value = 12

$ws.Range("D109").Value = 44511
$ws.Range("J109").Value = 100
$ws.Range("K109").Value = 3000
$ws.Range("L109").Value = 3000
$ws.Range("M109").Value = 3000
$ws.Range("N109").Value = "$/docena de atados (4 kilos)"
$ws.Range("O109").Value = "Región del Maule"
$ws.Range("P109").Value = 750
$ws.Range("Q109").Value = 4

$ws.Range("D110").Value = 44231
$ws.Range("J110").Value = 60
$ws.Range("K110").Value = 800
$ws.Range("L110").Value = 800
$ws.Range("M110").Value = 800
$ws.Range("N110").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O110").Value = "Región de La Araucanía"
$ws.Range("P110").Value = 533
$ws.Range("Q110").Value = 1.5

$ws.Range("D111").Value = 44400
$ws.Range("J111").Value = 300
$ws.Range("K111").Value = 1200
$ws.Range("L111").Value = 1200
$ws.Range("M111").Value = 1200
$ws.Range("N111").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O111").Value = "Región Metropolitana"
$ws.Range("P111").Value = 800
$ws.Range("Q111").Value = 1.5

$ws.Range("D112").Value = 44334
$ws.Range("J112").Value = 120
$ws.Range("K112").Value = 800
$ws.Range("L112").Value = 800
$ws.Range("M112").Value = 800
$ws.Range("N112").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O112").Value = "Región de La Araucanía"
$ws.Range("P112").Value = 533
$ws.Range("Q112").Value = 1.5

$ws.Range("D113").Value = 44573
$ws.Range("J113").Value = 20
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 10000
$ws.Range("N113").Value = "$/docena de atados (12 kilos)"
$ws.Range("O113").Value = "Región de La Araucanía"
$ws.Range("P113").Value = 833
$ws.Range("Q113").Value = 12

$ws.Range("D114").Value = 44403
$ws.Range("J114").Value = 150
$ws.Range("K114").Value = 4000
$ws.Range("L114").Value = 4000
$ws.Range("M114").Value = 4000
$ws.Range("N114").Value = "$/docena de atados (4 kilos)"
$ws.Range("O114").Value = "Región del Maule"
$ws.Range("P114").Value = 1000
$ws.Range("Q114").Value = 4

$ws.Range("D115").Value = 44319
$ws.Range("J115").Value = 60
$ws.Range("K115").Value = 800
$ws.Range("L115").Value = 800
$ws.Range("M115").Value = 800
$ws.Range("N115").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O115").Value = "Región de La Araucanía"
$ws.Range("P115").Value = 533
$ws.Range("Q115").Value = 1.5

$ws.Range("D116").Value = 44280
$ws.Range("J116").Value = 120
$ws.Range("K116").Value = 800
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 800
$ws.Range("N116").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O116").Value = "Región de La Araucanía"
$ws.Range("P116").Value = 533
$ws.Range("Q116").Value = 1.5

$ws.Range("D117").Value = 44362
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 3500
$ws.Range("L117").Value = 3500
$ws.Range("M117").Value = 3500
$ws.Range("N117").Value = "$/docena de atados (4 kilos)"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 875
$ws.Range("Q117").Value = 4

$ws.Range("D118").Value = 44365
$ws.Range("J118").Value = 200
$ws.Range("K118").Value = 3500
$ws.Range("L118").Value = 3500
$ws.Range("M118").Value = 3500
$ws.Range("N118").Value = "$/docena de atados (4 kilos)"
$ws.Range("O118").Value = "Región del Maule"
$ws.Range("P118").Value = 875
$ws.Range("Q118").Value = 4

$ws.Range("D119").Value = 44567
$ws.Range("J119").Value = 40
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 10000
$ws.Range("M119").Value = 10000
$ws.Range("N119").Value = "$/docena de atados (12 kilos)"
$ws.Range("O119").Value = "Región de La Araucanía"
$ws.Range("P119").Value = 833
$ws.Range("Q119").Value = 12

$ws.Range("D120").Value = 44473
$ws.Range("J120").Value = 100
$ws.Range("K120").Value = 4000
$ws.Range("L120").Value = 4000
$ws.Range("M120").Value = 4000
$ws.Range("N120").Value = "$/docena de atados (4 kilos)"
$ws.Range("O120").Value = "Región del Maule"
$ws.Range("P120").Value = 1000
$ws.Range("Q120").Value = 4

$ws.Range("D121").Value = 44537
$ws.Range("J121").Value = 200
$ws.Range("K121").Value = 3500
$ws.Range("L121").Value = 3500
$ws.Range("M121").Value = 3500
$ws.Range("N121").Value = "$/docena de atados (4 kilos)"
$ws.Range("O121").Value = "Región del Maule"
$ws.Range("P121").Value = 875
$ws.Range("Q121").Value = 4

$ws.Range("D122").Value = 44553
$ws.Range("J122").Value = 50
$ws.Range("K122").Value = 10000
$ws.Range("L122").Value = 10000
$ws.Range("M122").Value = 10000
$ws.Range("N122").Value = "$/docena de atados (12 kilos)"
$ws.Range("O122").Value = "Región de La Araucanía"
$ws.Range("P122").Value = 833
$ws.Range("Q122").Value = 12

$ws.Range("D123").Value = 44462
$ws.Range("J123").Value = 80
$ws.Range("K123").Value = 4000
$ws.Range("L123").Value = 4000
$ws.Range("M123").Value = 4000
$ws.Range("N123").Value = "$/docena de atados (4 kilos)"
$ws.Range("O123").Value = "Región del Maule"
$ws.Range("P123").Value = 1000
$ws.Range("Q123").Value = 4

$ws.Range("D124").Value = 44446
$ws.Range("J124").Value = 200
$ws.Range("K124").Value = 4000
$ws.Range("L124").Value = 4000
$ws.Range("M124").Value = 4000
$ws.Range("N124").Value = "$/docena de atados (4 kilos)"
$ws.Range("O124").Value = "Región del Maule"
$ws.Range("P124").Value = 1000
$ws.Range("Q124").Value = 4

$ws.Range("D125").Value = 44421
$ws.Range("J125").Value = 200
$ws.Range("K125").Value = 4000
$ws.Range("L125").Value = 4000
$ws.Range("M125").Value = 4000
$ws.Range("N125").Value = "$/docena de atados (4 kilos)"
$ws.Range("O125").Value = "Región del Maule"
$ws.Range("P125").Value = 1000
$ws.Range("Q125").Value = 4

$ws.Range("D126").Value = 44208
$ws.Range("J126").Value = 150
$ws.Range("K126").Value = 800
$ws.Range("L126").Value = 800
$ws.Range("M126").Value = 800
$ws.Range("N126").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O126").Value = "Región de La Araucanía"
$ws.Range("P126").Value = 533
$ws.Range("Q126").Value = 1.5

$ws.Range("D127").Value = 44355
$ws.Range("J127").Value = 290
$ws.Range("K127").Value = 800
$ws.Range("L127").Value = 800
$ws.Range("M127").Value = 800
$ws.Range("N127").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O127").Value = "Región de La Araucanía"
$ws.Range("P127").Value = 533
$ws.Range("Q127").Value = 1.5

$ws.Range("D128").Value = 44530
$ws.Range("J128").Value = 250
$ws.Range("K128").Value = 3500
$ws.Range("L128").Value = 3500
$ws.Range("M128").Value = 3500
$ws.Range("N128").Value = "$/docena de atados (4 kilos)"
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 875
$ws.Range("Q128").Value = 4

$ws.Range("D129").Value = 44483
$ws.Range("J129").Value = 100
$ws.Range("K129").Value = 4000
$ws.Range("L129").Value = 4000
$ws.Range("M129").Value = 4000
$ws.Range("N129").Value = "$/docena de atados (4 kilos)"
$ws.Range("O129").Value = "Región del Maule"
$ws.Range("P129").Value = 1000
$ws.Range("Q129").Value = 4

$ws.Range("D130").Value = 44617
$ws.Range("J130").Value = 120
$ws.Range("K130").Value = 10000
$ws.Range("L130").Value = 10000
$ws.Range("M130").Value = 10000
$ws.Range("N130").Value = "$/docena de atados (12 kilos)"
$ws.Range("O130").Value = "Región de La Araucanía"
$ws.Range("P130").Value = 833
$ws.Range("Q130").Value = 12

$ws.Range("D131").Value = 44264
$ws.Range("J131").Value = 150
$ws.Range("K131").Value = 800
$ws.Range("L131").Value = 800
$ws.Range("M131").Value = 800
$ws.Range("N131").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O131").Value = "Región de La Araucanía"
$ws.Range("P131").Value = 533
$ws.Range("Q131").Value = 1.5

$ws.Range("D132").Value = 44232
$ws.Range("J132").Value = 150
$ws.Range("K132").Value = 800
$ws.Range("L132").Value = 800
$ws.Range("M132").Value = 800
$ws.Range("N132").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O132").Value = "Región de La Araucanía"
$ws.Range("P132").Value = 533
$ws.Range("Q132").Value = 1.5

$ws.Range("D133").Value = 44279
$ws.Range("J133").Value = 50
$ws.Range("K133").Value = 800
$ws.Range("L133").Value = 800
$ws.Range("M133").Value = 800
$ws.Range("N133").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O133").Value = "Región de La Araucanía"
$ws.Range("P133").Value = 533
$ws.Range("Q133").Value = 1.5

$ws.Range("D134").Value = 44330
$ws.Range("J134").Value = 280
$ws.Range("K134").Value = 800
$ws.Range("L134").Value = 800
$ws.Range("M134").Value = 800
$ws.Range("N134").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O134").Value = "Región de La Araucanía"
$ws.Range("P134").Value = 533
$ws.Range("Q134").Value = 1.5

$ws.Range("D135").Value = 44504
$ws.Range("J135").Value = 100
$ws.Range("K135").Value = 3000
$ws.Range("L135").Value = 3000
$ws.Range("M135").Value = 3000
$ws.Range("N135").Value = "$/docena de atados (4 kilos)"
$ws.Range("O135").Value = "Región del Maule"
$ws.Range("P135").Value = 750
$ws.Range("Q135").Value = 4

$ws.Range("D136").Value = 44257
$ws.Range("J136").Value = 250
$ws.Range("K136").Value = 800
$ws.Range("L136").Value = 800
$ws.Range("M136").Value = 800
$ws.Range("N136").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O136").Value = "Región de La Araucanía"
$ws.Range("P136").Value = 533
$ws.Range("Q136").Value = 1.5

$ws.Range("D137").Value = 44301
$ws.Range("J137").Value = 84
$ws.Range("K137").Value = 800
$ws.Range("L137").Value = 800
$ws.Range("M137").Value = 800
$ws.Range("N137").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O137").Value = "Región de La Araucanía"
$ws.Range("P137").Value = 533
$ws.Range("Q137").Value = 1.5

$ws.Range("D138").Value = 44487
$ws.Range("J138").Value = 100
$ws.Range("K138").Value = 3500
$ws.Range("L138").Value = 3500
$ws.Range("M138").Value = 3500
$ws.Range("N138").Value = "$/docena de atados (4 kilos)"
$ws.Range("O138").Value = "Región del Maule"
$ws.Range("P138").Value = 875
$ws.Range("Q138").Value = 4

$ws.Range("D139").Value = 44236
$ws.Range("J139").Value = 250
$ws.Range("K139").Value = 800
$ws.Range("L139").Value = 800
$ws.Range("M139").Value = 800
$ws.Range("N139").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O139").Value = "Región de La Araucanía"
$ws.Range("P139").Value = 533
$ws.Range("Q139").Value = 1.5

$ws.Range("D140").Value = 44229
$ws.Range("J140").Value = 250
$ws.Range("K140").Value = 700
$ws.Range("L140").Value = 800
$ws.Range("M140").Value = 760
$ws.Range("N140").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O140").Value = "Provincia de Cautín"
$ws.Range("P140").Value = 507
$ws.Range("Q140").Value = 1.5

$ws.Range("D141").Value = 44299
$ws.Range("J141").Value = 300
$ws.Range("K141").Value = 800
$ws.Range("L141").Value = 800
$ws.Range("M141").Value = 800
$ws.Range("N141").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O141").Value = "Región de La Araucanía"
$ws.Range("P141").Value = 533
$ws.Range("Q141").Value = 1.5

$ws.Range("D142").Value = 44312
$ws.Range("J142").Value = 60
$ws.Range("K142").Value = 800
$ws.Range("L142").Value = 800
$ws.Range("M142").Value = 800
$ws.Range("N142").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O142").Value = "Región de La Araucanía"
$ws.Range("P142").Value = 533
$ws.Range("Q142").Value = 1.5

$ws.Range("D143").Value = 44399
$ws.Range("J143").Value = 200
$ws.Range("K143").Value = 3500
$ws.Range("L143").Value = 3500
$ws.Range("M143").Value = 3500
$ws.Range("N143").Value = "$/docena de atados (4 kilos)"
$ws.Range("O143").Value = "Región del Maule"
$ws.Range("P143").Value = 875
$ws.Range("Q143").Value = 4

$ws.Range("D144").Value = 44615
$ws.Range("J144").Value = 20
$ws.Range("K144").Value = 10000
$ws.Range("L144").Value = 10000
$ws.Range("M144").Value = 10000
$ws.Range("N144").Value = "$/docena de atados (12 kilos)"
$ws.Range("O144").Value = "Región de La Araucanía"
$ws.Range("P144").Value = 833
$ws.Range("Q144").Value = 12

$ws.Range("D145").Value = 44277
$ws.Range("J145").Value = 120
$ws.Range("K145").Value = 800
$ws.Range("L145").Value = 800
$ws.Range("M145").Value = 800
$ws.Range("N145").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O145").Value = "Región de La Araucanía"
$ws.Range("P145").Value = 533
$ws.Range("Q145").Value = 1.5

$ws.Range("D146").Value = 44390
$ws.Range("J146").Value = 250
$ws.Range("K146").Value = 3500
$ws.Range("L146").Value = 3500
$ws.Range("M146").Value = 3500
$ws.Range("N146").Value = "$/docena de atados (4 kilos)"
$ws.Range("O146").Value = "Región del Maule"
$ws.Range("P146").Value = 875
$ws.Range("Q146").Value = 4

$ws.Range("D147").Value = 44285
$ws.Range("J147").Value = 200
$ws.Range("K147").Value = 800
$ws.Range("L147").Value = 800
$ws.Range("M147").Value = 800
$ws.Range("N147").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O147").Value = "Región de La Araucanía"
$ws.Range("P147").Value = 533
$ws.Range("Q147").Value = 1.5

$ws.Range("D148").Value = 44498
$ws.Range("J148").Value = 200
$ws.Range("K148").Value = 3000
$ws.Range("L148").Value = 3000
$ws.Range("M148").Value = 3000
$ws.Range("N148").Value = "$/docena de atados (4 kilos)"
$ws.Range("O148").Value = "Región del Maule"
$ws.Range("P148").Value = 750
$ws.Range("Q148").Value = 4

$ws.Range("D149").Value = 44418
$ws.Range("J149").Value = 200
$ws.Range("K149").Value = 4000
$ws.Range("L149").Value = 4000
$ws.Range("M149").Value = 4000
$ws.Range("N149").Value = "$/docena de atados (4 kilos)"
$ws.Range("O149").Value = "Región del Maule"
$ws.Range("P149").Value = 1000
$ws.Range("Q149").Value = 4

$ws.Range("D150").Value = 44595
$ws.Range("J150").Value = 40
$ws.Range("K150").Value = 9000
$ws.Range("L150").Value = 9000
$ws.Range("M150").Value = 9000
$ws.Range("N150").Value = "$/docena de atados (12 kilos)"
$ws.Range("O150").Value = "Región de La Araucanía"
$ws.Range("P150").Value = 750
$ws.Range("Q150").Value = 12

$ws.Range("D151").Value = 44552
$ws.Range("J151").Value = 20
$ws.Range("K151").Value = 10000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 10000
$ws.Range("N151").Value = "$/docena de atados (12 kilos)"
$ws.Range("O151").Value = "Región de La Araucanía"
$ws.Range("P151").Value = 833
$ws.Range("Q151").Value = 12

$ws.Range("D152").Value = 44544
$ws.Range("J152").Value = 50
$ws.Range("K152").Value = 10000
$ws.Range("L152").Value = 10000
$ws.Range("M152").Value = 10000
$ws.Range("N152").Value = "$/docena de atados (12 kilos)"
$ws.Range("O152").Value = "Región de La Araucanía"
$ws.Range("P152").Value = 833
$ws.Range("Q152").Value = 12

